$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Finish existing row 68 (5:09PM debug entry): add Stop time + notes ---
$ws.Range("C68").NumberFormat = "h:mm"
$ws.Range("C68").Value = "5:12PM"
$ws.Range("H68").Value = "Quick fix. 2 variables not named correctly"

# --- New log entry in row 69 ---
$ws.Range("B69").Value = "5:12PM"
$ws.Range("C69").NumberFormat = "h:mm"
$ws.Range("C69").Value = 0.22916666666666666
$ws.Range("F69").Value = "Code"
$ws.Range("G69").Value = "Setting up code for newInvoiceCandS page functionality"

$ws.Range("H69").Select() | Out-Null
